# Update the "想去人数" (number of people interested) counts for two
# entries that appear on both the "展览" sheet and the "全部类型" sheet.
#   - Row 2 (南宁·布谷鸟动漫展5th):              F2  411 -> 412
#   - Row 3 (南宁·2024良牙动漫秋季盛典（秋典）): F3 5137 -> 5150

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 412
    $ws.Range("F3").Value = 5150
}
